$d = $word.ActiveDocument

# Step 1: remove the existing hidden _GoBack bookmark from the end of
# the "etcetera." paragraph; it will be re-added at the end of the
# last newly inserted paragraph below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Step 2: locate paragraph 7 ("... etcétera.") and append two brand
# new list-paragraphs right after it, inheriting its ListParagraph /
# numbering / justification / language formatting via InsertParagraphAfter,
# then replace that empty paragraph's content with the literal OOXML for
# each new bullet (this preserves <w:proofErr> spell-check wrappers and
# multi-run splits exactly as authored).
$srcPara = $d.Paragraphs.Item(7)

$null = $srcPara.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item(8)
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>Es posible desarrollar un programa que muestre vectores de tendencia con el SVM que ayuden a definir qué productos deberían salir al mercado</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> por medio del análisis de estos, así como la identificación de problemas más frecuentes y ver qué influye en la variación de las transacciones.</w:t></w:r></w:p>'
$newPara1.Range.InsertXML($xml1)

$newPara1 = $d.Paragraphs.Item(8)
$null = $newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item(9)
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">Se puede construir un </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">programa que aplique la técnica del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>bagging</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> para optimizar y mejorar un modelo de predicción de compra de los clientes actuales junto con métodos de ensamble para diferenciarlos y analizar las mejores opciones de venta.</w:t></w:r></w:p>'
$newPara2.Range.InsertXML($xml2)

# Step 3: the document's final paragraph is the old empty bold
# paragraph (w:rPr/w:b, no text). Re-use it in place for the third new
# bullet instead of appending + deleting, since it already sits at the
# end of the body (right before the sectPr) -- this also re-creates the
# _GoBack bookmark at the very end of the document's content.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>Otra opción es implementar un software que utilice los datos con el método de agrupamiento K-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>means</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> para determinar la influencia de cada grupo de datos </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>de acuerdo</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> a sus características, de esta forma se podría entrar en más detalle si se necesita más información y análisis</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> (teniendo en cuenta que este método es NP-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>hard</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$lastPara.Range.InsertXML($xml3)

Write-Host "Paragraphs now:" $d.Paragraphs.Count
